$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "224.92", "0.5310") are preserved exactly as text, matching the
# original inline-string cell contents, instead of being coerced to Double.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.351.50"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.716.59"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "224.92"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "0.5310"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("D9").Value = "0.2660"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").Value = "0.07697"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "4.495"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("D13").Value = "1.953.81"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "1.715.54"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "0.5822"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "0.0₅8227"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "67.99"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "27.386.43"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "222.96"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "4.671"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "6.030"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "145.05"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "1.707"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").Value = "7.258"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").Value = "16.26"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").Value = "0.05416"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "1.296"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "3.489"
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("D33").Value = "3.419"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "1.639"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "2.863"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "0.9558"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "2.394"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "0.5910"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").Value = "1.147.68"
$ws.Range("E39").Value = "  +9.46%  "
$ws.Range("D40").Value = "0.01657"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").Value = "5.835"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "0.8425"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").Value = "101.13"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "1.860.75"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "0.0₈117"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "57.99"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").Value = "0.4589"
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.185"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").Value = "1.008"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").Value = "0.05203"
$ws.Range("E51").Value = "  -0.89%  "

# Restore the default (unstyled) cell style now that the text values are set,
# so the cells end up with no explicit style reference, same as the source file.
$priceVolRange.Style = "Normal"
